# Update the payroll (nómina) row for employee 1009 / Gustavo Hernandez
# (row 6 of the "Reporte de Nómina" sheet) with the new figures from the
# latest frontend: ID Nómina, Horas Extras, Comisiones, Valor Horas Extras,
# Total Devengado, IGSS, Deducciones and Total Pagar all change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Nómina")

$ws.Range("H6").Value = "3011"
$ws.Range("K6").Value = "2"
$ws.Range("L6").Value = "Q500.00"
$ws.Range("N6").Value = "Q106.25"
$ws.Range("O6").Value = "Q9,356.25"
$ws.Range("T6").Value = "Q439.83"
$ws.Range("U6").Value = "Q439.83"
$ws.Range("V6").Value = "Q8,916.42"
